$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(318294931, Shalev  Afanasenko: 3,9)"
$ws.Range("B1").Value = "(305487936, Avihai  Kipnis: -4,2)"
$ws.Range("C1").Value = "(313227928, Aviv  Levi: 7,-8)"
$ws.Range("D1").Value = "(205807308, Sariel  Basis: -4,3)"
$ws.Range("E1").Value = "(315891549, Raz  Halaby: 8,-5)"
$ws.Range("F1").Value = "(315060103, Dan  Mshelh: -4,3)"
$ws.Range("G1").Value = "(313925141, Elad   Amer: -1,-5)"

$ws.Range("A3").Value = "cost: 657.4515679530903"
$ws.Range("A4").Value = "time: 90.35022399329858"
